# [FIX]: set new improve plan template.
# Update the two Velocity #foreach loops that render the "improve actions"
# table so they iterate over the new ${descriptions} / ${states} string
# collections instead of ${improveActions}, and collapse the now-redundant
# trailing "#end #end" down to a single "#end".
#
# NOTE: single-quoted PowerShell strings are used throughout so that the
# literal "${...}" Velocity placeholders are not expanded as PowerShell
# variables.

$d = $word.ActiveDocument

# 1) Description column: iterate ${descriptions}, print the bare ${item}.
$d.Content.Find.Execute(
    '#foreach(${item} in ${improveActions})${item.description} ',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    '#foreach( ${item} in ${descriptions} )${item} ',
    2) | Out-Null

# 2) Status column: iterate ${states}, print the bare ${item} (the old
#    #if/#else Completado/No Completada logic moves into ${states} itself).
$d.Content.Find.Execute(
    '#foreach(${item} in ${improveActions})#if(${item.completed} == true)Completado #else No Completada ',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    '#foreach( ${item} in ${states} )${item} ',
    2) | Out-Null

# 3) The status column used to close both the #if and the #foreach
#    ("#end #end"); now there is only the #foreach to close.
$d.Content.Find.Execute(
    '#end #end',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    '#end',
    2) | Out-Null
